# Insert a new row at position 17, shifting existing rows 17-124 down to 18-125.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(17).Insert()

# Fill in the new row 17 with the new weekly data entry.
# Columns A,B,C,E,F,G,H,I,Q,R mirror the rest of the table (same market/category template).
$ws.Range("A17").Value = 4
$ws.Range("B17").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C17").Value = 'Los Lagos'
$ws.Range("D17").Value = 44971
$ws.Range("E17").Value = 10
$ws.Range("F17").Value = 100112031
$ws.Range("G17").Value = 'Poroto verde'
$ws.Range("H17").Value = 'Magnum'
$ws.Range("I17").Value = 'Primera'
$ws.Range("J17").Value = 40
$ws.Range("K17").Value = 30000
$ws.Range("L17").Value = 30000
$ws.Range("M17").Value = 30000
$ws.Range("N17").Value = '$/saco 25 kilos'
$ws.Range("O17").Value = 'Región Metropolitana'
$ws.Range("P17").Value = 1200
$ws.Range("Q17").Value = 25
$ws.Range("R17").Value = 'Hortaliza'
